$d = $word.ActiveDocument

# The signature/date paragraph currently reads "{{SEDE}}, {{DATA_ASSINATURA}}".
# Drop the "{{SEDE}}, " portion so only the "{{DATA_ASSINATURA}}" placeholder remains.
$d.Content.Find.Execute("{{SEDE}}, {{DATA_ASSINATURA}}", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "{{DATA_ASSINATURA}}", 2)
